$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new rows 4 and 5, and a new "remark" column (E).
# Values are entered in an order that matches the first-occurrence order
# of the new shared strings in the target workbook:
#   0.0.8, 0.0.9, remark, "Fix array out of bound ...", version3

$ws.Range("A4").Value = "0.0.8"
$ws.Range("A5").Value = "0.0.9"
$ws.Range("E1").Value = "remark"
$ws.Range("E4").Value = "Fix array out of bound in ModifiedEnvironmentalSelection.java"
$ws.Range("C5").Value = "version3"

$ws.Range("B4").Value = "nsga-iii, m-nsga-iii"
$ws.Range("C4").Value = "version2"
$ws.Range("D4").Value = "version2"

$ws.Range("B5").Value = "nsga-iii, m-nsga-iii"
$ws.Range("D5").Value = "version3"

# Widen the new remark column (closest achievable width to 29.7109375 chars).
$ws.Columns.Item(5).ColumnWidth = 28.91

# Match the trailing selection left behind in the authored workbook.
$ws.Range("D8").Select()
